$d = $word.ActiveDocument

# 1) Update the existing sentence: r/s : Relationship -> e1/e2 : Entity wording.
$oldText = "anEntity.flatMap(ID::assert(r : Relationship) : s : Relationship (anEntity if equals, previous / next Entity if not equals)."
$newText = "anEntity.flatMap(ID::assert(e1 : Entity) : e2 : Entity (anEntity if same Entity, previous / next Entity if not same Entity)."

$rng = $d.Content
$found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

if (-not $found) {
    throw "Could not find target sentence to replace."
}

# 2) Insert two new paragraphs (each preceded by a blank paragraph) right
#    after that sentence's paragraph, before the following blank paragraph.
$rng = $d.Content
$found2 = $rng.Find.Execute($newText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)

if (-not $found2) {
    throw "Could not re-find updated sentence to anchor insertion."
}

$anchorIndex = $rng.Paragraphs(1).Index

$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.InsertParagraphAfter()
$rng.InsertParagraphAfter()
$rng.InsertParagraphAfter()

$d.Paragraphs($anchorIndex + 2).Range.Text = "Entity Monad built in Relationships: equals, inverseOf, parent, child, previous, next. Apply Relationship assert in the same manner."
$d.Paragraphs($anchorIndex + 4).Range.Text = "Domain Models Entities / Relationships: transforms of underlying Entities given Relationships contents."

Write-Output "done"
